$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5
$ws.Range("K2").Value = 4.6
$ws.Range("T3").Value = 1.75
$ws.Range("U3").Value = 2.22
$ws.Range("X3").Value = 21
$ws.Range("AC3").Value = 9.4
$ws.Range("AG3").Value = 21
$ws.Range("J4").Value = 3.8
$ws.Range("N4").Value = 5.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 1.62
$ws.Range("R4").Value = 1.61
$ws.Range("S4").Value = 2.54
$ws.Range("U4").Value = 2.66
$ws.Range("Z4").Value = 18
$ws.Range("AC4").Value = 9
$ws.Range("AD4").Value = 12
$ws.Range("AE4").Value = 22
$ws.Range("AG4").Value = 14.5
$ws.Range("AO4").Value = 12
$ws.Range("F5").Value = 2.02
$ws.Range("H5").Value = 1.6
$ws.Range("I5").Value = 4.9
$ws.Range("Q5").Value = 1.63
$ws.Range("F9").Value = 2.04
$ws.Range("I9").Value = 5.4
$ws.Range("J9").Value = 2.94
$ws.Range("P9").Value = 1.46
$ws.Range("Q9").Value = 2.78
$ws.Range("F10").Value = 2.86
$ws.Range("G10").Value = 3.3
$ws.Range("H10").Value = 2.9
$ws.Range("I10").Value = 3.5
$ws.Range("J10").Value = 2.6
$ws.Range("K10").Value = 3.2
$ws.Range("Q10").Value = 3
$ws.Range("F11").Value = 1.88
$ws.Range("G11").Value = 2.1
$ws.Range("H11").Value = 4
$ws.Range("I11").Value = 5.7
$ws.Range("J11").Value = 2.86
$ws.Range("K11").Value = 3.6
